# Update the "last updated" timestamp shown in A1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 8 de Septiembre de 2020 a las 15:22"

# Swap the display order of Banglades / Francia (Banglades overtook Francia
# in total cases). Row 17 now shows Banglades (with its updated stats),
# row 18 now shows Francia (kept at its previous stats / row position).
$ws.Cells.Item(17, 1).Value = "Banglades"
$ws.Cells.Item(18, 1).Value = "Francia"

# Estados Unidos (row 4)
$ws.Cells.Item(4, 2).Value = 6486426
$ws.Cells.Item(4, 3).Value = 851
$ws.Cells.Item(4, 4).Value = 3758629
$ws.Cells.Item(4, 5).Value = 2534211
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(4, 7).Value = 52
$ws.Cells.Item(4, 8).Value = 193586

# Banglades (row 17, new/updated data)
$ws.Cells.Item(17, 2).Value = 329251
$ws.Cells.Item(17, 3).Value = 1892
$ws.Cells.Item(17, 4).Value = 227809
$ws.Cells.Item(17, 5).Value = 96890
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 36
$ws.Cells.Item(17, 8).Value = 4552

# Francia (row 18, previous Banglades-row data, now shifted down)
$ws.Cells.Item(18, 2).Value = 328980
$ws.Cells.Item(18, 3).Value = 0
$ws.Cells.Item(18, 4).Value = 87836
$ws.Cells.Item(18, 5).Value = 210418
$ws.Cells.Item(18, 6).Value = 0
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 30726

# Arabia Saudita (row 19)
$ws.Cells.Item(19, 2).Value = 322237
$ws.Cells.Item(19, 3).Value = 781
$ws.Cells.Item(19, 4).Value = 298246
$ws.Cells.Item(19, 5).Value = 19854
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(19, 7).Value = 30
$ws.Cells.Item(19, 8).Value = 4137

# Oman (row 40)
$ws.Cells.Item(40, 2).Value = 87590
$ws.Cells.Item(40, 3).Value = 262
$ws.Cells.Item(40, 4).Value = 82973
$ws.Cells.Item(40, 5).Value = 3875
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 8
$ws.Cells.Item(40, 8).Value = 742

# Suecia (row 41)
$ws.Cells.Item(41, 2).Value = 85707
$ws.Cells.Item(41, 3).Value = 0
$ws.Cells.Item(41, 4).Value = 0
$ws.Cells.Item(41, 5).Value = 0
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(41, 7).Value = 4
$ws.Cells.Item(41, 8).Value = 5838

# Paises Bajos (row 44)
$ws.Cells.Item(44, 2).Value = 76548
$ws.Cells.Item(44, 3).Value = 964
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(44, 7).Value = 1
$ws.Cells.Item(44, 8).Value = 6244

# Suiza (row 62)
$ws.Cells.Item(62, 2).Value = 44837
$ws.Cells.Item(62, 3).Value = 245
$ws.Cells.Item(62, 4).Value = 37700
$ws.Cells.Item(62, 5).Value = 5121
$ws.Cells.Item(62, 6).Value = 0
$ws.Cells.Item(62, 7).Value = 2
$ws.Cells.Item(62, 8).Value = 2016

# Bosnia y Herzegovina (row 77)
$ws.Cells.Item(77, 2).Value = 21961
$ws.Cells.Item(77, 3).Value = 301
$ws.Cells.Item(77, 4).Value = 15172
$ws.Cells.Item(77, 5).Value = 6120
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(77, 7).Value = 5
$ws.Cells.Item(77, 8).Value = 669

# Dinamarca (row 83)
$ws.Cells.Item(83, 2).Value = 18356
$ws.Cells.Item(83, 3).Value = 243
$ws.Cells.Item(83, 4).Value = 15907
$ws.Cells.Item(83, 5).Value = 1821
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(83, 7).Value = 0
$ws.Cells.Item(83, 8).Value = 628

# Madagascar (row 85)
$ws.Cells.Item(85, 2).Value = 15435
$ws.Cells.Item(85, 3).Value = 83
$ws.Cells.Item(85, 4).Value = 14219
$ws.Cells.Item(85, 5).Value = 1010
$ws.Cells.Item(85, 6).Value = 0
$ws.Cells.Item(85, 7).Value = 4
$ws.Cells.Item(85, 8).Value = 206

# Republica de Macedonia (row 86)
$ws.Cells.Item(86, 2).Value = 15226
$ws.Cells.Item(86, 3).Value = 99
$ws.Cells.Item(86, 4).Value = 12700
$ws.Cells.Item(86, 5).Value = 1895
$ws.Cells.Item(86, 6).Value = 0
$ws.Cells.Item(86, 7).Value = 8
$ws.Cells.Item(86, 8).Value = 631

# Sri Lanka (row 133)
$ws.Cells.Item(133, 2).Value = 3126
$ws.Cells.Item(133, 3).Value = 3
$ws.Cells.Item(133, 4).Value = 2926
$ws.Cells.Item(133, 5).Value = 188
$ws.Cells.Item(133, 6).Value = 0
$ws.Cells.Item(133, 7).Value = 0
$ws.Cells.Item(133, 8).Value = 12

# Gibraltar (row 183)
$ws.Cells.Item(183, 2).Value = 320
$ws.Cells.Item(183, 3).Value = 5
$ws.Cells.Item(183, 4).Value = 276
$ws.Cells.Item(183, 5).Value = 44
$ws.Cells.Item(183, 6).Value = 0
$ws.Cells.Item(183, 7).Value = 0
$ws.Cells.Item(183, 8).Value = 0

# Bonaire, San Eustaquio y Saba (row 209)
$ws.Cells.Item(209, 2).Value = 21
$ws.Cells.Item(209, 3).Value = 3
$ws.Cells.Item(209, 4).Value = 7
$ws.Cells.Item(209, 5).Value = 14
$ws.Cells.Item(209, 6).Value = 0
$ws.Cells.Item(209, 7).Value = 0
$ws.Cells.Item(209, 8).Value = 0
